$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 7 data rows (old rows 2-8), shifting remaining rows up.
$ws.Rows("2:8").Delete()

# Append 6 new data rows at the end (now rows 16-21).
$newData = @(
    @(-2.083661317825317, -0.8663596510887146, -1.775479793548584),
    @(-4.610199928283691, 3.776979923248291, -1.897194743156433),
    @(0.9382890462875366, 0.2889392673969269, 0.8718574047088623),
    @(0.5700899958610535, -1.021672129631042, 0.4189008474349975),
    @(0.1930334120988845, -0.1259909570217132, 2.802496910095215),
    @(-0.2191479057073593, -2.771648406982422, 3.70489764213562)
)

$startRow = 16
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
